$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "1.003") must be
# explicitly formatted as Text first, otherwise Excel auto-converts the
# typed value into a real number and trailing zeros / exact text are lost.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D12", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D26", "D27", "D28", "D29", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.868.92'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.646.82'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '308.80'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").Value = '0.3886'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").Value = '0.3825'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").Value = '51.31'
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '0.08424'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").Value = '7.066'
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").Value = '7.890'
$ws.Range("E15").Value = '  +3.43%  '
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '1.650.93'
$ws.Range("E17").Value = '  +1.81%  '
$ws.Range("D18").Value = '94.28'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '0.06979'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = '19.63'
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '6.922'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("D24").Value = '23.886.93'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").Value = '2.965'
$ws.Range("E26").Value = '  +4.41%  '
$ws.Range("D27").Value = '22.00'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = '150.70'
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("D29").Value = '5.383'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '138.36'
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").Value = '2.508'
$ws.Range("E32").Value = '  +1.07%  '
$ws.Range("D33").Value = '1.830.84'
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("D34").Value = '1.051'
$ws.Range("E34").Value = '  +5.69%  '
$ws.Range("D35").Value = '0.08024'
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").Value = '0.02954'
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").Value = '6.715'
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").Value = '10.85'
$ws.Range("E38").Value = '  +4.85%  '
$ws.Range("D39").Value = '0.2677'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '0.09094'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '0.7560'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").Value = '13.44'
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").Value = '1.421'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").Value = '16.40'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").Value = '0.6937'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = '2.455'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").Value = '4.081'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").Value = '0.08266'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '133.91'
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").Value = '1.206'
$ws.Range("E51").Value = '  -0.43%  '
